$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 45-53 are the "2021-01-09" prediction-day block of weekly KNN forecasts
# (columns A, B, D, F only). Duplicate that block into rows 54-62 to extend
# the series, then append one brand-new week (row 63: "14 Mar -- 20 Mar 2021").
#
# Copy/paste (rather than literal .Value assignment) is used for column A so
# the "2021-01-09" text stays a shared string instead of being auto-converted
# into a date serial number by Excel's smart-typing, and for column F so the
# existing "KNN" shared string is reused without touching any other columns.

$ws.Range("A45:A53").Copy($ws.Range("A54:A62"))
$ws.Range("B45:B53").Copy($ws.Range("B54:B62"))
$ws.Range("D45:D53").Copy($ws.Range("D54:D62"))
$ws.Range("F45:F53").Copy($ws.Range("F54:F62"))

$ws.Cells.Item(53, 1).Copy($ws.Cells.Item(63, 1))
$ws.Cells.Item(63, 2).Value = "14 Mar -- 20 Mar 2021"
$ws.Cells.Item(63, 4).Value = 42.3
$ws.Cells.Item(53, 6).Copy($ws.Cells.Item(63, 6))
